$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.965.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.20%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.846.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.99%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.40%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("B5").Value = "'BNB"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'310.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.29%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "'USDC"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'1.012"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.37%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +2.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +2.18%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07231"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.26%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9289"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.68%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.34%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07738"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.82%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.805.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.35%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.344"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.45%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.427"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.83%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'88.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.15%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +0.44%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008638"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.83%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.40%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'27.000.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.29%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'5.059"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.66%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.79%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.917"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.50%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'152.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'18.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.42%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.000"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.14%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'114.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.39%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.967"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.87%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.08892"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.73%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.324"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +5.48%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.172"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.17%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7430"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.76%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.502"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.753"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.118"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.84%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.01960"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.36%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.05273"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.29%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.976"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.97%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.5214"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.95%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.985"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.45%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1512"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.86%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.93%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'10.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +6.08%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.4754"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.95%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.41%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'101.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.34%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.608"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.21%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'65.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.54%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06059"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.55%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.8880"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.13%  "
$ws.Range("E51").Style = "Normal"
